# Apply "add DQ Indicator for outliers" changes.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: DQ_Report
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("DQ_Report")

# Row 2
$ws1.Range("E2").Value = "Kodierung ist nicht eindeutig. Relation E75.0 - 846 ist im BfArM nicht vorhanden.  ICD10-Orpha Zuordnung ist gemäß BfArM nicht plausible. "

# Row 3
$ws1.Range("E3").Value = "Kodierung ist nicht eindeutig. Relation E75.0 - 797 ist im BfArM nicht vorhanden.  ICD10-Orpha Zuordnung ist gemäß BfArM nicht plausible. "

# Row 4
$ws1.Range("E4").Value = "Orpha Code 309151 ist im BfArM-Mapping nicht enthalten.  Kodierung ist nicht eindeutig. Relation E75.0 - 309151 ist im BfArM nicht vorhanden.  ICD10-Orpha Zuordnung ist gemäß BfArM nicht plausible. "

# Row 5
$ws1.Range("E5").Value = "Orpha Code 309247 ist im BfArM-Mapping nicht enthalten.  Kodierung ist nicht eindeutig. Relation E75.0 - 309247 ist im BfArM nicht vorhanden.  ICD10-Orpha Zuordnung ist gemäß BfArM nicht plausible. "

# Row 6
$ws1.Range("C6").Value = "G70"
$ws1.Range("E6").Value = "Kodierung ist nicht eindeutig. ICD10 Code G70 ist im BfArM Mapping nicht enthalten.  ICD10-Orpha Zuordnung ist gemäß BfArM nicht plausible. "

# Row 7
$ws1.Range("C7").Value = "G70"
$ws1.Range("E7").Value = "Kodierung ist nicht eindeutig. ICD10 Code G70 ist im BfArM Mapping nicht enthalten.  ICD10-Orpha Zuordnung ist gemäß BfArM nicht plausible. "

# Row 8
$ws1.Range("C8").Value = "E84.80"
$ws1.Range("D8").Value = 589
$ws1.Range("E8").Value = "Kodierung ist nicht eindeutig. Relation E84.80 - 589 ist im BfArM nicht vorhanden.  ICD10-Orpha Zuordnung ist gemäß BfArM nicht plausible. "

# Row 10
$ws1.Range("E10").Value = "Orpha Code 3 ist im BfArM-Mapping nicht enthalten.  Fehlendes ICD10 Code.  "

# Row 12
$ws1.Range("E12").Value = "Orpha Code 320 ist im BfArM-Mapping nicht enthalten.  Kodierung ist nicht eindeutig. Relation E66.89 - 320 ist im BfArM nicht vorhanden.  ICD10-Orpha Zuordnung ist gemäß BfArM nicht plausible. "

# Row 13
$ws1.Range("E13").Value = "Kodierung ist nicht eindeutig. Relation G35.9 - 71529 ist im BfArM nicht vorhanden.  ICD10-Orpha Zuordnung ist gemäß BfArM nicht plausible. "

# Row 17
$ws1.Range("E17").Value = "Kodierung ist nicht eindeutig. Relation E75.2 - 342 ist im BfArM nicht vorhanden.  ICD10-Orpha Zuordnung ist gemäß BfArM nicht plausible. "

# Row 18
$ws1.Range("E18").Value = "Kodierung ist nicht eindeutig. Relation E75.0 - 226 ist im BfArM nicht vorhanden.  ICD10-Orpha Zuordnung ist gemäß BfArM nicht plausible. "

# ---------------------------------------------------------------------------
# Sheet 2: Statistik - new "outlier" / orphaCoding indicator columns
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Statistik")

# Shift the last two header/value columns (pt_no, case_no) two slots to the
# right (H->J, I->K) before inserting the new columns in their place, so no
# data is lost.
$ws2.Range("K1").Value = "case_no"
$ws2.Range("K2").Value = 28
$ws2.Range("J1").Value = "pt_no"
$ws2.Range("J2").Value = 27

# Header row updates / additions
$ws2.Range("C1").Value = "outlier_rate"
$ws2.Range("F1").Value = "orphaCoding_plausibility_rate"
$ws2.Range("G1").Value = "orphaCoding_relativeFrequency"
$ws2.Range("H1").Value = "orphaCoding_absoluteFrequency"
$ws2.Range("I1").Value = "uniqueRd_no"

# Data row updates / additions
$ws2.Range("B2").Value = 40.82
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 66.67
$ws2.Range("F2").Value = 31.25
$ws2.Range("G2").Value = 1.4
$ws2.Range("H2").Value = 14
$ws2.Range("I2").Value = 9
